$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and 1h volume change (E) values
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.225.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.982.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.428'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.364'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.490.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000162'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.224.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.981.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.84'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '329.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.492'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.47'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.104.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0919'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.76'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.75%  '
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0658'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.016.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.168.33'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('E46').Value = '  -3.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.918'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0235'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.38'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0848'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.55%  '
